{"js": "// Mill and Ref License Format fix\n// Applies the trader-license edits: new company name, new address,\n// new trader type description, trims the \"withdraw purchased sugar...\"\n// sentence, and updates the ST license number and TIN.\n\nasync function replaceOnce(context, searchText, replacementText, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const results = context.document.body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacementText, \"Replace\");\n  await context.sync();\n}\n\n// 1. Company name in the title line.\nawait replaceOnce(\n  context,\n  \"CENTRAL AZUCARERA DON PEDRO, INC.\",\n  \"DAFFID INDUSTRIES, INC.\"\n);\n\n// 2. Street address.\nawait replaceOnce(\n  context,\n  \" 14/F, Net One Center, 3rd Ave. cor. 26th St., E. Square, CPW, BGC, Taguig\",\n  \" 802 Atlanta Center, #31 Annapolis St., Greenhills, SanJuan City\"\n);\n\n// 3. Trader type / commodity description.\nawait replaceOnce(\n  context,\n  \", is hereby licensed with this Office to operate as a DOMESTIC SUGAR TRADER during the \",\n  \", is hereby licensed with this Office to operate as INTERNATIONAL SUGAR TRADER for Chemically Pure Fructose and High Fructose Corn Syrup during the \"\n);\n\n// 4. Drop the \"Said Trader is hereby authorized to withdraw purchased sugar\n//    from the warehouse...\" sentence, leaving only \" Crop Year.\" This span\n//    crosses three runs (two of which are fully removed); search() can\n//    match across run boundaries, and insertText(\"\", \"Replace\") deletes it.\nawait replaceOnce(\n  context,\n  \" Said Trader is hereby authorized to withdraw purchased sugar from the warehouse of any mill or refinery subject to rules and regulations issued by this Office pursuant thereto.\",\n  \"\"\n);\n\n// 5. ST license number.\nawait replaceOnce(context, \"   ST-2021-002\", \"   ST-2021-001-F\");\n\n// 6. TIN number.\nawait replaceOnce(context, \"214-280-422-000\", \"234-562-022-000\");\n", "ps1": "# Mill and Ref License Format fix\n# Applies the trader-license edits: new company name, new address,\n# new trader type description, trims the \"withdraw purchased sugar...\"\n# sentence, and updates the ST license number and TIN.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once {\n    param([string]$FindText, [string]$ReplaceText)\n    # wdFindContinue = 1, wdReplaceAll = 2 (numeric literals; this host has\n    # no predefined wdXxx constant variables).\n    $found = $d.Content.Find.Execute(\n        $FindText,\n        $false, $false, $false, $false, $false,\n        $true, 1, $false,\n        $ReplaceText,\n        2)\n    if (-not $found) {\n        throw \"Text not found: $FindText\"\n    }\n}\n\n# 1. Company name in the title line.\nReplace-Once \"CENTRAL AZUCARERA DON PEDRO, INC.\" \"DAFFID INDUSTRIES, INC.\"\n\n# 2. Street address.\nReplace-Once \" 14/F, Net One Center, 3rd Ave. cor. 26th St., E. Square, CPW, BGC, Taguig\" \" 802 Atlanta Center, #31 Annapolis St., Greenhills, SanJuan City\"\n\n# 3. Trader type / commodity description.\nReplace-Once \", is hereby licensed with this Office to operate as a DOMESTIC SUGAR TRADER during the \" \", is hereby licensed with this Office to operate as INTERNATIONAL SUGAR TRADER for Chemically Pure Fructose and High Fructose Corn Syrup during the \"\n\n# 4. Drop the \"Said Trader is hereby authorized to withdraw purchased sugar\n#    from the warehouse...\" sentence, leaving only \" Crop Year.\" This span\n#    crosses three runs (two of which are fully removed); Find.Execute can\n#    match across run boundaries, and an empty ReplaceWith deletes it.\nReplace-Once \" Said Trader is hereby authorized to withdraw purchased sugar from the warehouse of any mill or refinery subject to rules and regulations issued by this Office pursuant thereto.\" \"\"\n\n# 5. ST license number.\nReplace-Once \"   ST-2021-002\" \"   ST-2021-001-F\"\n\n# 6. TIN number.\nReplace-Once \"214-280-422-000\" \"234-562-022-000\"\n"}
